$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.314.16"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.53%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.834.12"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").Value = "  +0.79%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4738"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3688"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07444"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.13%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8857"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.93%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.46"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.871.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.61%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07342"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.35%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.441"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "93.85"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.01%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.575"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.10%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008791"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.566.69"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.81"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.289"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +1.07%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.090.94"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.21%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.893"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "152.14"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +1.44%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.151"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.224"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.67%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "117.17"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08992"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7499"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.176"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.550"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("E37").Value = "  +1.23%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05346"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.34%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01956"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.969"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.392"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "7.222"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.62%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5309"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  +0.67%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4932"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "105.23"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +0.94%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06302"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
